$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 584
$ws.Range("A584").Value = "21TRD09200"
$ws.Range("B584").Value = "Bunner"
$ws.Range("C584").Value = "DUS UCM"
$ws.Range("D584").NumberFormat = "@"
$ws.Range("D584").Value = "4510.111"
$ws.Range("E584").Value = "UCM"
$ws.Range("F584").Value = "Guilty"
$ws.Range("G584").Value = "Guilty"
$ws.Range("H584").NumberFormat = "@"
$ws.Range("H584").Value = "`$ 0"
$ws.Range("I584").NumberFormat = "@"
$ws.Range("I584").Value = "`$ 0"

# Row 585
$ws.Range("A585").Value = "21TRD09200"
$ws.Range("B585").Value = "Bunner"
$ws.Range("C585").Value = "OPERATING W/O A VALID OL - UCM"
$ws.Range("D585").NumberFormat = "@"
$ws.Range("D585").Value = "4510.12"
$ws.Range("E585").Value = "UCM"
$ws.Range("F585").Value = "Guilty"
$ws.Range("G585").Value = "Guilty"
$ws.Range("H585").NumberFormat = "@"
$ws.Range("H585").Value = "`$ 0"
$ws.Range("I585").NumberFormat = "@"
$ws.Range("I585").Value = "`$ 0"

# Row 586
$ws.Range("A586").Value = "21TRD09200"
$ws.Range("B586").Value = "Bunner"
$ws.Range("C586").Value = "FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS"
$ws.Range("D586").Value = "4510.21A*"
$ws.Range("E586").Value = "UCM"
$ws.Range("F586").Value = "Guilty"
$ws.Range("G586").Value = "Guilty"
$ws.Range("H586").NumberFormat = "@"
$ws.Range("H586").Value = "`$ 0"
$ws.Range("I586").NumberFormat = "@"
$ws.Range("I586").Value = "`$ 0"

# Row 587
$ws.Range("A587").Value = "21TRD09200"
$ws.Range("B587").Value = "Bunner"
$ws.Range("C587").Value = "FAILURE TO FILE REGISTRATION"
$ws.Range("D587").NumberFormat = "@"
$ws.Range("D587").Value = "4503.11"
$ws.Range("E587").Value = "MM"
$ws.Range("F587").Value = "Guilty"
$ws.Range("G587").Value = "Guilty"
$ws.Range("H587").NumberFormat = "@"
$ws.Range("H587").Value = "`$ 0"
$ws.Range("I587").NumberFormat = "@"
$ws.Range("I587").Value = "`$ 0"

# Row 588
$ws.Range("A588").Value = "21TRD09200"
$ws.Range("B588").Value = "Bunner"
$ws.Range("C588").Value = "DUS UCM"
$ws.Range("D588").NumberFormat = "@"
$ws.Range("D588").Value = "4510.111"
$ws.Range("E588").Value = "UCM"
$ws.Range("F588").Value = "Guilty"
$ws.Range("G588").Value = "Guilty"
$ws.Range("H588").NumberFormat = "@"
$ws.Range("H588").Value = "`$ 0"
$ws.Range("I588").NumberFormat = "@"
$ws.Range("I588").Value = "`$ 0"
$ws.Range("J588").Value = "None"
$ws.Range("K588").Value = "None"

# Row 589
$ws.Range("A589").Value = "21TRD09200"
$ws.Range("B589").Value = "Bunner"
$ws.Range("C589").Value = "OPERATING W/O A VALID OL - UCM"
$ws.Range("D589").NumberFormat = "@"
$ws.Range("D589").Value = "4510.12"
$ws.Range("E589").Value = "UCM"
$ws.Range("F589").Value = "Guilty"
$ws.Range("G589").Value = "Guilty"
$ws.Range("H589").NumberFormat = "@"
$ws.Range("H589").Value = "`$ 0"
$ws.Range("I589").NumberFormat = "@"
$ws.Range("I589").Value = "`$ 0"
$ws.Range("J589").Value = "None"
$ws.Range("K589").Value = "None"

# Row 590
$ws.Range("A590").Value = "21TRD09200"
$ws.Range("B590").Value = "Bunner"
$ws.Range("C590").Value = "FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS"
$ws.Range("D590").Value = "4510.21A*"
$ws.Range("E590").Value = "UCM"
$ws.Range("F590").Value = "Guilty"
$ws.Range("G590").Value = "Guilty"
$ws.Range("H590").NumberFormat = "@"
$ws.Range("H590").Value = "`$ 0"
$ws.Range("I590").NumberFormat = "@"
$ws.Range("I590").Value = "`$ 0"
$ws.Range("J590").Value = "None"
$ws.Range("K590").Value = "None"

# Row 591
$ws.Range("A591").Value = "21TRD09200"
$ws.Range("B591").Value = "Bunner"
$ws.Range("C591").Value = "FAILURE TO FILE REGISTRATION"
$ws.Range("D591").NumberFormat = "@"
$ws.Range("D591").Value = "4503.11"
$ws.Range("E591").Value = "MM"
$ws.Range("F591").Value = "Guilty"
$ws.Range("G591").Value = "Guilty"
$ws.Range("H591").NumberFormat = "@"
$ws.Range("H591").Value = "`$ 0"
$ws.Range("I591").NumberFormat = "@"
$ws.Range("I591").Value = "`$ 0"
$ws.Range("J591").Value = "None"
$ws.Range("K591").Value = "None"

# Row 592
$ws.Range("A592").Value = "21TRD09200"
$ws.Range("B592").Value = "Bunner"
$ws.Range("C592").Value = "DUS UCM"
$ws.Range("D592").NumberFormat = "@"
$ws.Range("D592").Value = "4510.111"
$ws.Range("E592").Value = "UCM"
$ws.Range("F592").Value = "Guilty"

# Row 593
$ws.Range("A593").Value = "21TRD09200"
$ws.Range("B593").Value = "Bunner"
$ws.Range("C593").Value = "OPERATING W/O A VALID OL - UCM"
$ws.Range("D593").NumberFormat = "@"
$ws.Range("D593").Value = "4510.12"
$ws.Range("E593").Value = "UCM"
$ws.Range("F593").Value = "Guilty"

# Row 594
$ws.Range("A594").Value = "21TRD09200"
$ws.Range("B594").Value = "Bunner"
$ws.Range("C594").Value = "FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS"
$ws.Range("D594").Value = "4510.21A*"
$ws.Range("E594").Value = "UCM"
$ws.Range("F594").Value = "Guilty"

# Row 595
$ws.Range("A595").Value = "21TRD09200"
$ws.Range("B595").Value = "Bunner"
$ws.Range("C595").Value = "FAILURE TO FILE REGISTRATION"
$ws.Range("D595").NumberFormat = "@"
$ws.Range("D595").Value = "4503.11"
$ws.Range("E595").Value = "MM"
$ws.Range("F595").Value = "Guilty"

# Row 596
$ws.Range("A596").Value = "21TRD09200"
$ws.Range("B596").Value = "Bunner"
$ws.Range("C596").Value = "DUS UCM"
$ws.Range("D596").NumberFormat = "@"
$ws.Range("D596").Value = "4510.111"
$ws.Range("E596").Value = "UCM"
$ws.Range("F596").Value = "Guilty"

# Row 597
$ws.Range("A597").Value = "21TRD09200"
$ws.Range("B597").Value = "Bunner"
$ws.Range("C597").Value = "OPERATING W/O A VALID OL - UCM"
$ws.Range("D597").NumberFormat = "@"
$ws.Range("D597").Value = "4510.12"
$ws.Range("E597").Value = "UCM"
$ws.Range("F597").Value = "Guilty"

# Row 598
$ws.Range("A598").Value = "21TRD09200"
$ws.Range("B598").Value = "Bunner"
$ws.Range("C598").Value = "FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS"
$ws.Range("D598").Value = "4510.21A*"
$ws.Range("E598").Value = "UCM"
$ws.Range("F598").Value = "Guilty"

# Row 599
$ws.Range("A599").Value = "21TRD09200"
$ws.Range("B599").Value = "Bunner"
$ws.Range("C599").Value = "FAILURE TO FILE REGISTRATION"
$ws.Range("D599").NumberFormat = "@"
$ws.Range("D599").Value = "4503.11"
$ws.Range("E599").Value = "MM"
$ws.Range("F599").Value = "Guilty"

# Row 600
$ws.Range("A600").Value = "21TRD09200"
$ws.Range("B600").Value = "Bunner"
$ws.Range("C600").Value = "DUS UCM"
$ws.Range("D600").NumberFormat = "@"
$ws.Range("D600").Value = "4510.111"
$ws.Range("E600").Value = "UCM"
$ws.Range("F600").Value = "No Contest"
$ws.Range("G600").Value = "Guilty"
$ws.Range("H600").NumberFormat = "@"
$ws.Range("H600").Value = "`$ 0"
$ws.Range("I600").NumberFormat = "@"
$ws.Range("I600").Value = "`$ 0"

# Row 601
$ws.Range("A601").Value = "21TRD09200"
$ws.Range("B601").Value = "Bunner"
$ws.Range("C601").Value = "OPERATING W/O A VALID OL - UCM"
$ws.Range("D601").NumberFormat = "@"
$ws.Range("D601").Value = "4510.12"
$ws.Range("E601").Value = "UCM"
$ws.Range("F601").Value = "No Contest"
$ws.Range("G601").Value = "Guilty"
$ws.Range("H601").NumberFormat = "@"
$ws.Range("H601").Value = "`$ 0"
$ws.Range("I601").NumberFormat = "@"
$ws.Range("I601").Value = "`$ 0"

# Row 602
$ws.Range("A602").Value = "21TRD09200"
$ws.Range("B602").Value = "Bunner"
$ws.Range("C602").Value = "FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS"
$ws.Range("D602").Value = "4510.21A*"
$ws.Range("E602").Value = "UCM"
$ws.Range("F602").Value = "No Contest"
$ws.Range("G602").Value = "Guilty"
$ws.Range("H602").NumberFormat = "@"
$ws.Range("H602").Value = "`$ 0"
$ws.Range("I602").NumberFormat = "@"
$ws.Range("I602").Value = "`$ 0"

# Row 603
$ws.Range("A603").Value = "21TRD09200"
$ws.Range("B603").Value = "Bunner"
$ws.Range("C603").Value = "FAILURE TO FILE REGISTRATION"
$ws.Range("D603").NumberFormat = "@"
$ws.Range("D603").Value = "4503.11"
$ws.Range("E603").Value = "MM"
$ws.Range("F603").Value = "No Contest"
$ws.Range("G603").Value = "Guilty"
$ws.Range("H603").NumberFormat = "@"
$ws.Range("H603").Value = "`$ 0"
$ws.Range("I603").NumberFormat = "@"
$ws.Range("I603").Value = "`$ 0"

# Row 604
$ws.Range("A604").Value = "21TRD09200"
$ws.Range("B604").Value = "Bunner"
$ws.Range("C604").Value = "DUS UCM"
$ws.Range("D604").NumberFormat = "@"
$ws.Range("D604").Value = "4510.111"
$ws.Range("E604").Value = "UCM"
$ws.Range("F604").Value = "No Contest"
$ws.Range("G604").Value = "Guilty"
$ws.Range("H604").NumberFormat = "@"
$ws.Range("H604").Value = "`$ 0"
$ws.Range("I604").NumberFormat = "@"
$ws.Range("I604").Value = "`$ 0"
$ws.Range("J604").Value = "None"
$ws.Range("K604").Value = "None"

# Row 605
$ws.Range("A605").Value = "21TRD09200"
$ws.Range("B605").Value = "Bunner"
$ws.Range("C605").Value = "OPERATING W/O A VALID OL - UCM"
$ws.Range("D605").NumberFormat = "@"
$ws.Range("D605").Value = "4510.12"
$ws.Range("E605").Value = "UCM"
$ws.Range("F605").Value = "No Contest"
$ws.Range("G605").Value = "Guilty"
$ws.Range("H605").NumberFormat = "@"
$ws.Range("H605").Value = "`$ 0"
$ws.Range("I605").NumberFormat = "@"
$ws.Range("I605").Value = "`$ 0"
$ws.Range("J605").Value = "None"
$ws.Range("K605").Value = "None"

# Row 606
$ws.Range("A606").Value = "21TRD09200"
$ws.Range("B606").Value = "Bunner"
$ws.Range("C606").Value = "FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS"
$ws.Range("D606").Value = "4510.21A*"
$ws.Range("E606").Value = "UCM"
$ws.Range("F606").Value = "No Contest"
$ws.Range("G606").Value = "Guilty"
$ws.Range("H606").NumberFormat = "@"
$ws.Range("H606").Value = "`$ 0"
$ws.Range("I606").NumberFormat = "@"
$ws.Range("I606").Value = "`$ 0"
$ws.Range("J606").Value = "None"
$ws.Range("K606").Value = "None"

# Row 607
$ws.Range("A607").Value = "21TRD09200"
$ws.Range("B607").Value = "Bunner"
$ws.Range("C607").Value = "FAILURE TO FILE REGISTRATION"
$ws.Range("D607").NumberFormat = "@"
$ws.Range("D607").Value = "4503.11"
$ws.Range("E607").Value = "MM"
$ws.Range("F607").Value = "No Contest"
$ws.Range("G607").Value = "Guilty"
$ws.Range("H607").NumberFormat = "@"
$ws.Range("H607").Value = "`$ 0"
$ws.Range("I607").NumberFormat = "@"
$ws.Range("I607").Value = "`$ 0"
$ws.Range("J607").Value = "None"
$ws.Range("K607").Value = "None"
